# Updates Price (D) and Volume(1h) (E) columns for the cryptos table
# to match the latest scraped snapshot.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "85.985.28"
$ws.Range("E2").Value = "  +7.06%  "
$ws.Range("D3").Value = "3.314.97"
$ws.Range("E3").Value = "  +3.26%  "
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "219.06"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +4.20%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "633.72"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.08%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.325"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +17.73%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.999"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.05%  "
$ws.Range("E9").Value = "  -1.71%  "
$ws.Range("D10").Value = "3.313.59"
$ws.Range("E10").Value = "  +3.27%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.596"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -3.70%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000276"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +5.19%  "
$ws.Range("E13").Value = "  +0.22%  "
$ws.Range("D14").Value = "3.925.40"
$ws.Range("E14").Value = "  +3.19%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "34.03"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +4.73%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.39"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.18%  "
$ws.Range("D17").Value = "85.481.14"
$ws.Range("E17").Value = "  +6.35%  "
$ws.Range("D18").Value = "3.308.76"
$ws.Range("E18").Value = "  +2.80%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "14.61"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.24%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "3.16"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +5.30%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "443.21"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.04%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.10"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.49%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.24"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.06%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "7.38"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +8.35%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "5.45"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +13.84%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "12.23"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +11.61%  "
$ws.Range("D27").Value = "3.485.56"
$ws.Range("E27").Value = "  +2.95%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "78.17"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.76%  "
$ws.Range("E29").Value = "  +3.83%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.999"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.01%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.168"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +33.92%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "608.70"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +9.30%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "9.22"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.56%  "
$ws.Range("E34").Value = "  +0.28%  "
$ws.Range("E35").Value = "  +2.76%  "
$ws.Range("E36").Value = "  +0.23%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.150"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.09%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "23.36"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.46%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "6.43"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +11.80%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.00"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.04%  "
$ws.Range("E41").Value = "  -0.90%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "21.27"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +3.18%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.05"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +11.58%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.07"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +12.40%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "158.47"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -4.45%  "
$ws.Range("E46").Value = "  +0.00%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "188.47"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.01%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.36"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.75%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "45.21"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +3.74%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.789"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.98%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "26.29"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +3.56%  "
